$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.907.86"
$ws.Range("E2").Value = "  +0.08%  "

$ws.Range("D3").Value = "1.636.83"
$ws.Range("E3").Value = "  +0.17%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.15%  "

$ws.Range("E6").Value = "  -0.43%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.40"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.91%  "

$ws.Range("E9").Value = "  -0.28%  "

$ws.Range("E10").Value = "  -0.03%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0883"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.22%  "

$ws.Range("D12").Value = "1.869.76"
$ws.Range("E12").Value = "  +0.19%  "

$ws.Range("D13").Value = "1.638.36"
$ws.Range("E13").Value = "  +0.83%  "

$ws.Range("E14").Value = "  -0.80%  "

$ws.Range("E15").Value = "  -0.92%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.38"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.38%  "

$ws.Range("D17").Value = "27.913.42"
$ws.Range("E17").Value = "  +0.07%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "228.97"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.20%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.70"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.75%  "

$ws.Range("D20").Value = "0.0₃0720"
$ws.Range("E20").Value = "  +0.06%  "

$ws.Range("E21").Value = "  +0.01%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.35"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.25%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.06"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.41%  "

$ws.Range("E24").Value = "  -0.02%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "155.69"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.81%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.89"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.91%  "

$ws.Range("E27").Value = "  +0.31%  "

$ws.Range("E28").Value = "  -0.35%  "

$ws.Range("E29").Value = "  -0.04%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.19"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.23%  "

$ws.Range("E31").Value = "  -0.03%  "

$ws.Range("E33").Value = "  +1.71%  "

$ws.Range("D34").Value = "1.396.61"
$ws.Range("E34").Value = "  -0.58%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.61"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.66%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.01"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.49%  "

$ws.Range("E37").Value = "  -0.73%  "

$ws.Range("E38").Value = "  +0.73%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.562"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.23%  "

$ws.Range("E40").Value = "  -2.28%  "

$ws.Range("E41").Value = "  -0.04%  "

$ws.Range("E42").Value = "  -1.24%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.84"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.31%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "65.98"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.25%  "

$ws.Range("E45").Value = "  -0.85%  "

$ws.Range("D46").Value = "1.776.06"
$ws.Range("E46").Value = "  -0.02%  "

$ws.Range("E47").Value = "  -2.51%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "88.70"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.21%  "

$ws.Range("E49").Value = "  +2.61%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0505"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.02%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.65"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.34%  "
